# Refresh the "想去人数" (want-to-go count) figures pulled from the source site.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) rows 3-20, column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 14
$wsExpo.Range("F4").Value  = 1365
$wsExpo.Range("F5").Value  = 305
$wsExpo.Range("F6").Value  = 1030
$wsExpo.Range("F7").Value  = 10636
$wsExpo.Range("F8").Value  = 18
$wsExpo.Range("F10").Value = 288
$wsExpo.Range("F12").Value = 697
$wsExpo.Range("F13").Value = 12028
$wsExpo.Range("F14").Value = 12455
$wsExpo.Range("F16").Value = 118
$wsExpo.Range("F19").Value = 74
$wsExpo.Range("F20").Value = 42

# Sheet "全部类型" (all types) rows 4-21, column F -- same underlying entries,
# shifted down by one row relative to "展览" because this sheet carries one
# extra leading row.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 14
$wsAll.Range("F5").Value  = 1365
$wsAll.Range("F6").Value  = 305
$wsAll.Range("F7").Value  = 1030
$wsAll.Range("F8").Value  = 10636
$wsAll.Range("F9").Value  = 18
$wsAll.Range("F11").Value = 288
$wsAll.Range("F13").Value = 697
$wsAll.Range("F14").Value = 12028
$wsAll.Range("F15").Value = 12455
$wsAll.Range("F17").Value = 118
$wsAll.Range("F20").Value = 74
$wsAll.Range("F21").Value = 42
